$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MAPPING")

# The ENTITY column (B) previously pointed every table at the single
# "CD_HISTORICAL" database. Split it out so feed/catalog tables point at
# "CD_RAW" and the rest continue to resolve through the renamed
# "CD_ADMIN" database.
$ws.Range("B2").Value = "CD_ADMIN"   # FRIENDLY
$ws.Range("B3").Value = "CD_ADMIN"   # HIGH_DATE_TABLES
$ws.Range("B4").Value = "CD_ADMIN"   # SLA_CONFIG
$ws.Range("B5").Value = "CD_ADMIN"   # RUN_HISTORY_TABLES
$ws.Range("B6").Value = "CD_RAW"     # DATA_FEED_CATALOG
$ws.Range("B7").Value = "CD_RAW"     # DATA_FEED_MAPPING

# Leave the MAPPING sheet as the active tab (it was RUN_HISTORY_TABLES
# before this edit).
$ws.Activate()
